# "Fruta / hortaliza, semanal" -- insert a new weekly price record for
# Albahaca (Terminal La Palmera de La Serena) as row 25, pushing the
# existing rows 25-35 down to 26-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 25; this shifts rows 25:35 -> 26:36.
$ws.Rows(25).Insert()

# Populate the newly-inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 8
$ws.Range("B25").Value = "Terminal La Palmera de La Serena"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 44455
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 100112052
$ws.Range("G25").Value = "Albahaca"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 600
$ws.Range("K25").Value = 4500
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 4750
$ws.Range("N25").Value = "`$/paquete"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 4750
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = "Hortaliza"
